$d = $word.ActiveDocument

# 1. Update the SDK version text
$d.Content.Find.Execute(
    "Windows 10 Creators Update SDK (15063)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Windows 10 Fall Creators Update SDK (16299)", 2)
